$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Within each 4-row year block (A, B, C, D periods), the "B" period row and
# "C" period row have their data (columns A-E) swapped. The "A" and "D"
# period rows are left untouched.
$pairs = @(
    @(3,4), @(7,8), @(11,12), @(15,16), @(19,20), @(23,24), @(27,28),
    @(31,32), @(35,36), @(39,40), @(43,44), @(47,48), @(51,52), @(55,56),
    @(59,60), @(63,64), @(67,68)
)

foreach ($pair in $pairs) {
    $rowB = $pair[0]
    $rowC = $pair[1]

    # Columns A, B, C, E always carry real content and can be swapped safely.
    $rngB = $ws.Range("A$($rowB):C$($rowB)")
    $rngC = $ws.Range("A$($rowC):C$($rowC)")
    $valsB = $rngB.Value2
    $valsC = $rngC.Value2
    $rngB.Value = $valsC
    $rngC.Value = $valsB

    $eB = $ws.Range("E$($rowB)")
    $eC = $ws.Range("E$($rowC)")
    $eValB = $eB.Value2
    $eValC = $eC.Value2
    $eB.Value = $eValC
    $eC.Value = $eValB

    # Column D only holds real numbers for the later (2016+) blocks; for the
    # earlier blocks it is an empty placeholder cell on both rows, so
    # swapping it is a no-op that we skip (writing "" would clear the cell
    # and change its stored type).
    $dB = $ws.Range("D$($rowB)")
    $dC = $ws.Range("D$($rowC)")
    $dValB = $dB.Value2
    $dValC = $dC.Value2
    if (($dValB -ne "") -or ($dValC -ne "")) {
        $dB.Value = $dValC
        $dC.Value = $dValB
    }
}

# Columns F (民用钢质船舶产销率) and G (民用钢质船舶销售量) are removed entirely.
$ws.Columns.Item(7).Delete()
$ws.Columns.Item(6).Delete()
